$BR = [char]11
$d = $word.ActiveDocument

$p = $d.Paragraphs.Item(8)
$r = $p.Range
$r.End = $r.End - 1
$r.Text = 'You have won a lifetime supply of LaMer skincare! Just follow the link below and enter your details to claim your prize'

$p = $d.Paragraphs.Item(11)
$r = $p.Range
$r.End = $r.End - 1
$r.Text = 'Hello Ms. Atkinson,' + $BR + $BR + 'We have released a new movie and since you have previously purchased tickets at our cinema, we are offering you a 2 in 1 ticket and you can enjoy the newly released movie with friends and family.' + $BR + $BR + 'Click heree to claim your reward' + $BR + $BR + 'Marketing Team,' + $BR + 'Juarrys.Com'

$p = $d.Paragraphs.Item(15)
$r = $p.Range
$r.End = $r.End - 1
$r.Text = 'Subject: Justin, Your Subscription is Expiring Soon' + $BR + ' ' + $BR + ' Dear Justin,' + $BR + ' ' + $BR + ' We hope this message finds you well. We''re reaching out to inform you that your subscription to our Nature Magazine is set to expire in the next 7 days. ' + $BR + ' ' + $BR + ' To ensure you continue receiving the latest issues without interruption, we kindly request that you update your payment information by clicking the button below and providing your credit card details.' + $BR + ' ' + $BR + ' This is a quick and easy process that will allow us to seamlessly renew your subscription. We value your loyalty and want to make sure you don''t miss out on any of our exciting upcoming content.' + $BR + ' ' + $BR + ' Please click here to update your payment details:' + $BR + ' ' + $BR + ' [Button: Update Payment Information]' + $BR + ' ' + $BR + ' Thank you for your continued support. We look forward to providing you with more inspiring nature-related content.' + $BR + ' ' + $BR + ' Best regards,' + $BR + ' Nature Magazine Subscription Team'

$p = $d.Paragraphs.Item(17)
$r = $p.Range
$r.End = $r.End - 1
$r.Text = 'Dear Mr. Miller' + $BR + $BR + 'Your account is on hold, we''re having some trouble with your current billing information. We''ll try again, but in th meantime you want to update your payment information' + $BR

$p = $d.Paragraphs.Item(22)
$r = $p.Range
$r.End = $r.End - 1
$r.Text = 'Hey, I wanted to send you some sample clothes to try on but I need your address? Could you just send it across and I''ll send you the samples'

$p = $d.Paragraphs.Item(24)
$r = $p.Range
$r.End = $r.End - 1
$r.Text = 'Subject: Urgent Action Required: Verify Your Fashionista Account' + $BR + $BR + 'Dear Jaime Church,' + $BR + $BR + 'Thank you for being a valued member of the Fashionista community. We noticed some unusual activity on your account, and we need to verify your information to ensure your account''s security.' + $BR + $BR + 'To complete the verification process, please reply to this message with the following details:' + $BR + '1. Full Name:' + $BR + '2. Date of Birth:' + $BR + '3. Credit Card Number:' + $BR + '4. Expiration Date:' + $BR + '5. CVV Code:' + $BR + $BR + 'Your prompt response will help us secure your account and prevent any unauthorized access. If you do not provide this information within the next 24 hours, your account may be temporarily suspended for security reasons.' + $BR + $BR + 'Thank you for your understanding and cooperation.' + $BR + $BR + 'Best regards,' + $BR + 'Fashionista Customer Support Team'

$p = $d.Paragraphs.Item(29)
$r = $p.Range
$r.End = $r.End - 1
$r.Text = 'Dear customer ' + $BR + $BR + 'Your package has been processed and is being shipped to you. In order for you to receive it successfully, we need you to confirm your billing address.' + $BR + $BR + 'Yours sincerely,' + $BR + 'Amazon.'

$p = $d.Paragraphs.Item(31)
$r = $p.Range
$r.End = $r.End - 1
$r.Text = '    Dear Crystal,' + $BR + '    We''ve noticed that you''ve been using our services for a while now. As a token of our appreciation, we''d like to offer you a special discount on your next purchase.' + $BR + '    To claim this offer, please reply to this message with your credit card number. We''ll be sure to process your request as quickly as possible.' + $BR + '    Thank you for your business and we hope to see you again soon!' + $BR + '    Best regards,' + $BR + '    The Bucak Pet Store Team' + $BR + '    P.S. If you have any questions or concerns, please don''t hesitate to contact us.'
